$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Linear and Quadratic coefficients driving the attenuation model
$ws.Range("M3").Value = 0.42
$ws.Range("N3").Value = 0

# Recalculate so cached formula results (and, where supported, chart
# caches) reflect the new inputs
$excel.CalculateFullRebuild()
$excel.Calculate()

# Best-effort nudge of the scatter chart's cached plot values (the
# "Brightness" series, Sheet1!$I$2:$I$27) so they pick up the recalculated
# numbers too
try {
    $chart = $ws.ChartObjects(1).Chart
    $chart.Refresh()
} catch {
}

# Match the new active cell selection recorded in the sheet view
$ws.Activate()
$ws.Range("M4").Select()
